# Applies the Golem_Profits profit-recalculation update across all class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 669.1429000000001
$ws.Range("I19").Value = 424.75
$ws.Range("J19").Value = 995
$ws.Range("K19").Value = 424.75
$ws.Range("L19").Value = 995
$ws.Range("M19").Value = -249.75
$ws.Range("N19").Value = -1345

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 311
$ws.Range("I31").Value = 284.7143
$ws.Range("J31").Value = 495
$ws.Range("K31").Value = 854.1428999999999
$ws.Range("L31").Value = 1485
$ws.Range("M31").Value = -624.1428999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3989
$ws.Range("I34").Value = 3989
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3989
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3786

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 3989
$ws.Range("I36").Value = 3989
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3989
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3274

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 899.25
$ws.Range("I43").Value = 898
$ws.Range("J43").Value = 899.6667
$ws.Range("K43").Value = 898
$ws.Range("L43").Value = 899.6667
$ws.Range("M43").Value = -829
$ws.Range("N43").Value = -1037.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3745

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3745

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2969
$ws.Range("I92").Value = 2157.6
$ws.Range("J92").Value = 4997.5
$ws.Range("K92").Value = 2157.6
$ws.Range("L92").Value = 4997.5
$ws.Range("M92").Value = -909.5999999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 561.6667
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 592.5
$ws.Range("K98").Value = 500
$ws.Range("L98").Value = 592.5
$ws.Range("M98").Value = 998
$ws.Range("N98").Value = -3588.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 561.6667
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 592.5
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 1777.5
$ws.Range("M122").Value = 950
$ws.Range("N122").Value = -6677.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1373.5
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 747
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 6723
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -11643

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5086.875
$ws.Range("I138").Value = 5678.5
$ws.Range("J138").Value = 3312
$ws.Range("K138").Value = 17035.5
$ws.Range("L138").Value = 9936
$ws.Range("M138").Value = -11895.5
$ws.Range("N138").Value = -20216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 39999
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 39999
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 39999
$ws.Range("N37").Value = -40545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2924.75
$ws.Range("I45").Value = 2924.75
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2924.75
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2547.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 15533.167
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 15533.167
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 15533.167
$ws.Range("N95").Value = -21025.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 24398.6
$ws.Range("I124").Value = 9998
$ws.Range("J124").Value = 45999.5
$ws.Range("K124").Value = 9998
$ws.Range("L124").Value = 45999.5
$ws.Range("M124").Value = -5088
$ws.Range("N124").Value = -55819.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 36000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 36000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 36000
$ws.Range("N88").Value = -36812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 36000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 36000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 36000
$ws.Range("N91").Value = -38808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.08
$ws.Range("I7").Value = 75.3125
$ws.Range("J7").Value = 155.22223
$ws.Range("K7").Value = 75.3125
$ws.Range("L7").Value = 155.22223
$ws.Range("M7").Value = 37.6875
$ws.Range("N7").Value = -381.22223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1499.5
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -649
$ws.Range("N22").Value = -2700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1669283.4
$ws.Range("I99").Value = 1252675
$ws.Range("J99").Value = 2502500
$ws.Range("K99").Value = 1252675
$ws.Range("L99").Value = 2502500
$ws.Range("M99").Value = -1251177
$ws.Range("N99").Value = -2505496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4068
$ws.Range("I122").Value = 547.3333
$ws.Range("J122").Value = 6708.5
$ws.Range("K122").Value = 1641.9999
$ws.Range("L122").Value = 20125.5
$ws.Range("M122").Value = 808.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1669283.4
$ws.Range("I126").Value = 1252675
$ws.Range("J126").Value = 2502500
$ws.Range("K126").Value = 3758025
$ws.Range("L126").Value = 7507500
$ws.Range("M126").Value = -3755555
$ws.Range("N126").Value = -7512440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 9514.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 9514.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 28543.5
$ws.Range("N105").Value = -33785.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 1030
$ws.Range("I125").Value = 1030
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3090
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 1830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3986.353
$ws.Range("I122").Value = 3025.25
$ws.Range("J122").Value = 6293
$ws.Range("K122").Value = 9075.75
$ws.Range("L122").Value = 18879
$ws.Range("M122").Value = -6625.75
$ws.Range("N122").Value = -23779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1178.8572
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 1300.3334
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 1300.3334
$ws.Range("M22").Value = -155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1178.8572
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 1300.3334
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 1300.3334
$ws.Range("M27").Value = -343

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3793.2666
$ws.Range("I122").Value = 3571.2856
$ws.Range("J122").Value = 3987.5
$ws.Range("K122").Value = 10713.8568
$ws.Range("L122").Value = 11962.5
$ws.Range("M122").Value = -8263.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 35476.332
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 35476.332
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 35476.332
$ws.Range("N124").Value = -45296.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1849.5
$ws.Range("I62").Value = 1199
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 1199
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -575
$ws.Range("N62").Value = -3748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 1849.5
$ws.Range("I65").Value = 1199
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 5995
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -2875
$ws.Range("N65").Value = -18740
